$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 79240.60000000001
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 79240.60000000001
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 79240.60000000001
$ws.Cells.Item(3, 14).Value = -79468.60000000001

$ws.Cells.Item(12, 8).Value = 201.5
$ws.Cells.Item(12, 9).Value = 159.38461
$ws.Cells.Item(12, 10).Value = 384
$ws.Cells.Item(12, 11).Value = 159.38461
$ws.Cells.Item(12, 12).Value = 384
$ws.Cells.Item(12, 13).Value = 10.61538999999999
$ws.Cells.Item(12, 14).Value = -724

$ws.Cells.Item(17, 8).Value = 853.32654
$ws.Cells.Item(17, 9).Value = 964.5
$ws.Cells.Item(17, 10).Value = 848.59576
$ws.Cells.Item(17, 11).Value = 2893.5
$ws.Cells.Item(17, 12).Value = 2545.78728
$ws.Cells.Item(17, 13).Value = -2725.5
$ws.Cells.Item(17, 14).Value = -2881.78728

$ws.Cells.Item(58, 8).Value = 1224.75
$ws.Cells.Item(58, 9).Value = 713.5
$ws.Cells.Item(58, 10).Value = 2758.5
$ws.Cells.Item(58, 11).Value = 2140.5
$ws.Cells.Item(58, 12).Value = 8275.5
$ws.Cells.Item(58, 13).Value = -1990.5
$ws.Cells.Item(58, 14).Value = -8575.5

$ws.Cells.Item(102, 8).Value = 79240.60000000001
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = 79240.60000000001
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 79240.60000000001
$ws.Cells.Item(102, 14).Value = -85730.60000000001

$ws.Cells.Item(112, 8).Value = 9373.393
$ws.Cells.Item(112, 9).Value = 645.8333
$ws.Cells.Item(112, 10).Value = 11753.637
$ws.Cells.Item(112, 11).Value = 1937.4999
$ws.Cells.Item(112, 12).Value = 35260.911
$ws.Cells.Item(112, 13).Value = -829.4999
$ws.Cells.Item(112, 14).Value = -37476.911

$ws.Cells.Item(137, 8).Value = 1120382.5
$ws.Cells.Item(137, 9).Value = 894.5625
$ws.Cells.Item(137, 10).Value = 2498213.8
$ws.Cells.Item(137, 11).Value = 2683.6875
$ws.Cells.Item(137, 12).Value = 7494641.399999999
$ws.Cells.Item(137, 13).Value = -133.6875
$ws.Cells.Item(137, 14).Value = -7499741.399999999

$ws.Cells.Item(138, 8).Value = 1838.8823
$ws.Cells.Item(138, 9).Value = 1342.871
$ws.Cells.Item(138, 10).Value = 2607.7
$ws.Cells.Item(138, 11).Value = 4028.613
$ws.Cells.Item(138, 12).Value = 7823.099999999999
$ws.Cells.Item(138, 13).Value = 1111.387
$ws.Cells.Item(138, 14).Value = -18103.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(57, 8).Value = 11121444
$ws.Cells.Item(57, 9).Value = 11121444
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 11121444
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = -11120960

$ws.Cells.Item(61, 8).Value = 1335490.4
$ws.Cells.Item(61, 9).Value = 1390844.1
$ws.Cells.Item(61, 10).Value = 7000
$ws.Cells.Item(61, 11).Value = 1390844.1
$ws.Cells.Item(61, 12).Value = 7000
$ws.Cells.Item(61, 13).Value = -1390632.1
$ws.Cells.Item(61, 14).Value = -7424

$ws.Cells.Item(74, 8).Value = 3050915.2
$ws.Cells.Item(74, 9).Value = 3572237.8
$ws.Cells.Item(74, 10).Value = 9866.666999999999
$ws.Cells.Item(74, 11).Value = 3572237.8
$ws.Cells.Item(74, 12).Value = 9866.666999999999
$ws.Cells.Item(74, 13).Value = -3571363.8
$ws.Cells.Item(74, 14).Value = -11614.667

$ws.Cells.Item(77, 8).Value = 3050915.2
$ws.Cells.Item(77, 9).Value = 3572237.8
$ws.Cells.Item(77, 10).Value = 9866.666999999999
$ws.Cells.Item(77, 11).Value = 17861189
$ws.Cells.Item(77, 12).Value = 49333.335
$ws.Cells.Item(77, 13).Value = -17856821
$ws.Cells.Item(77, 14).Value = -58069.335

$ws.Cells.Item(97, 8).Value = 1989.4
$ws.Cells.Item(97, 9).Value = 1989.4
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1989.4
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -1493.4

$ws.Cells.Item(122, 8).Value = 2888.1843
$ws.Cells.Item(122, 9).Value = 2494.8215
$ws.Cells.Item(122, 10).Value = 3989.6
$ws.Cells.Item(122, 11).Value = 7484.4645
$ws.Cells.Item(122, 12).Value = 11968.8
$ws.Cells.Item(122, 13).Value = -5034.4645
$ws.Cells.Item(122, 14).Value = -16868.8

$ws.Cells.Item(124, 8).Value = 30000
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 30000
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 30000
$ws.Cells.Item(124, 13).Value = ""
$ws.Cells.Item(124, 14).Value = -39820

$ws.Cells.Item(132, 8).Value = 1189305.5
$ws.Cells.Item(132, 9).Value = 1826294.8
$ws.Cells.Item(132, 10).Value = 6325.5713
$ws.Cells.Item(132, 11).Value = 5478884.4
$ws.Cells.Item(132, 12).Value = 18976.7139
$ws.Cells.Item(132, 13).Value = -5476354.4
$ws.Cells.Item(132, 14).Value = -24036.7139

$ws.Cells.Item(134, 8).Value = 79000
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 79000
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 79000
$ws.Cells.Item(134, 14).Value = -89140

$ws.Cells.Item(136, 8).Value = 1335490.4
$ws.Cells.Item(136, 9).Value = 1390844.1
$ws.Cells.Item(136, 10).Value = 7000
$ws.Cells.Item(136, 11).Value = 4172532.3
$ws.Cells.Item(136, 12).Value = 21000
$ws.Cells.Item(136, 13).Value = -4169982.3
$ws.Cells.Item(136, 14).Value = -26100

$ws.Cells.Item(139, 8).Value = 99639
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 99639
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 99639
$ws.Cells.Item(139, 14).Value = -109919

$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(13, 8).Value = 73424.5
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 73424.5
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 73424.5
$ws.Cells.Item(13, 14).Value = -73760.5

$ws.Cells.Item(22, 8).Value = 350.75
$ws.Cells.Item(22, 9).Value = 350.75
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 350.75
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -177.75

$ws.Cells.Item(118, 8).Value = 45355.5
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 45355.5
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 45355.5
$ws.Cells.Item(118, 14).Value = -48669.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 2166.6667
$ws.Cells.Item(3, 9).Value = 2000
$ws.Cells.Item(3, 10).Value = 2500
$ws.Cells.Item(3, 11).Value = 2000
$ws.Cells.Item(3, 12).Value = 2500
$ws.Cells.Item(3, 13).Value = -1887
$ws.Cells.Item(3, 14).Value = -2726

$ws.Cells.Item(58, 8).Value = 652806.5600000001
$ws.Cells.Item(58, 9).Value = 1123593.6
$ws.Cells.Item(58, 10).Value = 5474.375
$ws.Cells.Item(58, 11).Value = 1123593.6
$ws.Cells.Item(58, 12).Value = 5474.375
$ws.Cells.Item(58, 13).Value = -1123390.6
$ws.Cells.Item(58, 14).Value = -5880.375

$ws.Cells.Item(132, 8).Value = 22855208
$ws.Cells.Item(132, 9).Value = 30316348
$ws.Cells.Item(132, 10).Value = 471790
$ws.Cells.Item(132, 11).Value = 90949044
$ws.Cells.Item(132, 12).Value = 1415370
$ws.Cells.Item(132, 13).Value = -90946514
$ws.Cells.Item(132, 14).Value = -1420430

$ws.Cells.Item(134, 8).Value = 29834.416
$ws.Cells.Item(134, 9).Value = 38644.777
$ws.Cells.Item(134, 10).Value = 3403.3333
$ws.Cells.Item(134, 11).Value = 115934.331
$ws.Cells.Item(134, 12).Value = 10209.9999
$ws.Cells.Item(134, 13).Value = -113399.331
$ws.Cells.Item(134, 14).Value = -15279.9999

$ws.Cells.Item(136, 8).Value = 652806.5600000001
$ws.Cells.Item(136, 9).Value = 1123593.6
$ws.Cells.Item(136, 10).Value = 5474.375
$ws.Cells.Item(136, 11).Value = 3370780.8
$ws.Cells.Item(136, 12).Value = 16423.125
$ws.Cells.Item(136, 13).Value = -3368230.8
$ws.Cells.Item(136, 14).Value = -21523.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 5
$ws.Cells.Item(12, 9).Value = 5
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(12, 11).Value = 15
$ws.Cells.Item(12, 12).Value = 15
$ws.Cells.Item(12, 13).Value = 158
$ws.Cells.Item(12, 14).Value = -361

$ws.Cells.Item(60, 8).Value = 516.4545000000001
$ws.Cells.Item(60, 9).Value = 203.22223
$ws.Cells.Item(60, 10).Value = 1926
$ws.Cells.Item(60, 11).Value = 609.66669
$ws.Cells.Item(60, 12).Value = 5778
$ws.Cells.Item(60, 13).Value = -358.66669
$ws.Cells.Item(60, 14).Value = -6280

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 37912.375
$ws.Cells.Item(10, 9).Value = 2699.75
$ws.Cells.Item(10, 10).Value = 73125
$ws.Cells.Item(10, 11).Value = 2699.75
$ws.Cells.Item(10, 12).Value = 73125
$ws.Cells.Item(10, 13).Value = -2530.75
$ws.Cells.Item(10, 14).Value = -73463

$ws.Cells.Item(102, 8).Value = 3998.125
$ws.Cells.Item(102, 9).Value = 3553.5278
$ws.Cells.Item(102, 10).Value = 7999.5
$ws.Cells.Item(102, 11).Value = 3553.5278
$ws.Cells.Item(102, 12).Value = 7999.5
$ws.Cells.Item(102, 13).Value = -1931.5278
$ws.Cells.Item(102, 14).Value = -11243.5

$ws.Cells.Item(123, 8).Value = 55902
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 55902
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 55902
$ws.Cells.Item(123, 14).Value = -60802

$ws.Cells.Item(133, 8).Value = 75142.25
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 75142.25
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 75142.25
$ws.Cells.Item(133, 14).Value = -85262.25

$ws.Cells.Item(136, 8).Value = 27415.625
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 27415.625
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 82246.875
$ws.Cells.Item(136, 14).Value = -87346.875

$ws.Cells.Item(141, 8).Value = 46905.6
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 46905.6
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 46905.6
$ws.Cells.Item(141, 14).Value = -57265.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 1062500
$ws.Cells.Item(20, 9).Value = 50000
$ws.Cells.Item(20, 10).Value = 1400000
$ws.Cells.Item(20, 11).Value = 50000
$ws.Cells.Item(20, 12).Value = 1400000
$ws.Cells.Item(20, 13).Value = -49774
$ws.Cells.Item(20, 14).Value = -1400452

$ws.Cells.Item(42, 8).Value = 18016.666
$ws.Cells.Item(42, 9).Value = 18016.666
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 18016.666
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -17453.666

$ws.Cells.Item(49, 8).Value = 18016.666
$ws.Cells.Item(49, 9).Value = 18016.666
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 18016.666
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = -17869.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 150698
$ws.Cells.Item(119, 9).Value = 0
$ws.Cells.Item(119, 10).Value = 150698
$ws.Cells.Item(119, 11).Value = 0
$ws.Cells.Item(119, 12).Value = 150698
$ws.Cells.Item(119, 14).Value = -160374

$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).Value = ""

$ws.Cells.Item(136, 8).Value = 1451935.8
$ws.Cells.Item(136, 9).Value = 1643194.1
$ws.Cells.Item(136, 10).Value = 49374
$ws.Cells.Item(136, 11).Value = 4929582.300000001
$ws.Cells.Item(136, 12).Value = 148122
$ws.Cells.Item(136, 13).Value = -4927032.300000001
$ws.Cells.Item(136, 14).Value = -153222

$ws.Cells.Item(140, 8).Value = 66429
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 66429
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 66429
$ws.Cells.Item(140, 14).Value = -76789

$ws.Cells.Item(141, 8).Value = 60715
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 60715
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 60715
$ws.Cells.Item(141, 14).Value = -71075

Write-Output "applied all changes"
